# Auto-generated edit script
# Updates "想去人数" (column F) values on sheets "展览" (1) and "全部类型" (4)

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet "展览" (Worksheets(1)) updates ---
$wsExhibit.Range("F3").Value = 503
$wsExhibit.Range("F4").Value = 1462
$wsExhibit.Range("F5").Value = 143
$wsExhibit.Range("F7").Value = 295
$wsExhibit.Range("F8").Value = 5210
$wsExhibit.Range("F9").Value = 125
$wsExhibit.Range("F10").Value = 705
$wsExhibit.Range("F12").Value = 53
$wsExhibit.Range("F13").Value = 297
$wsExhibit.Range("F14").Value = 41
$wsExhibit.Range("F15").Value = 6256
$wsExhibit.Range("F17").Value = 135
$wsExhibit.Range("F18").Value = 137
$wsExhibit.Range("F20").Value = 15049
$wsExhibit.Range("F21").Value = 1494
$wsExhibit.Range("F22").Value = 262
$wsExhibit.Range("F25").Value = 10950
$wsExhibit.Range("F26").Value = 716
$wsExhibit.Range("F27").Value = 4277
$wsExhibit.Range("F28").Value = 212
$wsExhibit.Range("F30").Value = 122

# --- Sheet "全部类型" (Worksheets(4)) updates ---
$wsAll.Range("F3").Value = 503
$wsAll.Range("F4").Value = 1462
$wsAll.Range("F5").Value = 143
$wsAll.Range("F7").Value = 295
$wsAll.Range("F9").Value = 5210
$wsAll.Range("F10").Value = 125
$wsAll.Range("F11").Value = 705
$wsAll.Range("F14").Value = 53
$wsAll.Range("F15").Value = 297
$wsAll.Range("F16").Value = 41
$wsAll.Range("F18").Value = 6256
$wsAll.Range("F20").Value = 135
$wsAll.Range("F21").Value = 137
$wsAll.Range("F23").Value = 15049
$wsAll.Range("F24").Value = 1494
$wsAll.Range("F25").Value = 262
$wsAll.Range("F28").Value = 10950
$wsAll.Range("F29").Value = 716
$wsAll.Range("F30").Value = 0
$wsAll.Range("F31").Value = 212
$wsAll.Range("F33").Value = 122

